$d = $word.ActiveDocument

function Replace-ParagraphText {
    param($doc, [string]$oldText, [string]$newText)
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs($i).Range
        $r = $doc.Range($para.Start, $para.End - 1)
        if ($r.Text -eq $oldText) {
            $r.Text = $newText
            Write-Host "Updated paragraph $i"
            return $true
        }
    }
    Write-Host "WARNING: no paragraph matched expected text: $($oldText.Substring(0, [Math]::Min(60, $oldText.Length)))"
    return $false
}

# Change 1
$old0 = "Strategic storytelling executive with 15+ years of experience shaping how enterprise B2B organizations communicate, perform, and grow—operating at the intersection of brand, GTM alignment, and customer experience to translate complex positioning into executive narratives that resonate with senior stakeholders. Umberto partners with C-suite and product leaders to elevate brand from a communications tool to a revenue-driving asset, building corporate messaging platforms that connect emotionally and scale commercially across complex sales cycles. His work spans competitive positioning, narrative development, and go-to-market strategies for global consulting firms, SaaS leaders, and financial services institutions—recognized by Cannes Lions, The Webby Awards, and AdAge. Known for cross-functional leadership that drives alignment across teams, he translates brand vision into enterprise-scale systems through strategic counsel and collaboration with executive-level audiences. His approach ensures marketing delivers measurable business impact—connecting insight with execution to drive both customer value and business growth."
$new0 = "Strategic storytelling executive with 15+ years of experience helping enterprise B2B organizations translate complex capabilities into competitive positioning that drives GTM alignment and revenue growth. Umberto partners with senior stakeholders and C-suite leaders to develop executive narratives and brand platforms that connect insight with execution—elevating brand from a communications tool to a strategic asset that supports complex sales cycles and enterprise buyer journeys. His work spans global consulting firms, SaaS leaders, financial services institutions, and iconic brands including Johnson & Johnson, American Express, and General Electric, with recognition from Cannes Lions, The Webby Awards, and AdAge. Known for cross-functional leadership and the ability to lead without formal authority, he builds corporate messaging systems and strategic communications infrastructure that align product, sales, and leadership teams around unified market narratives. From executive-level presentations to campaign architecture, his approach ensures every touchpoint delivers measurable business impact while scaling commercially across global markets."
Replace-ParagraphText $d $old0 $new0 | Out-Null

# Change 2
$old1 = "Created differentiated, category-distinct brand strategy and executive narrative platform to reposition Synovus as a relationship-driven B2B banking partner amid declining market share and mindshare. Led messaging, identity refresh, and agency-led strategic storytelling to drive GTM alignment and brand clarity—lifting brand favorability by 12%, boosting digital engagement by 11%, achieving 9% increase in new customer acquisition, and contributing to 64.3% revenue growth between 2021 and 2023."
$new1 = "Led strategic storytelling and brand repositioning for B2B banking partner facing declining market share, developing category-distinct creative platform and executive narratives that drove GTM alignment across messaging, identity, and agency-led campaigns. The new platform drove a 9% increase in new customer acquisition, lifted brand favorability by 12%, boosted digital engagement by 11%, and contributed to 64.3% revenue growth between 2021 and 2023."
Replace-ParagraphText $d $old1 $new1 | Out-Null

# Change 3
$old2 = "Transformed CRM strategy across financial services, hospitality, and retail portfolios by replacing legacy journey-based models with a dynamic, moments-based engagement framework—translating AI-powered predictive personalization into compelling narratives that drove stakeholder adoption. Designed and implemented automated content systems and GenAI-enabled tools that accelerated execution, improved message relevance, and unlocked personalized brand experiences at scale through strategic storytelling that connected technical capabilities to business outcomes. Delivered 5% lift in customer retention, 12% increase in cross-sell performance, and 32% reduction in production costs—demonstrating measurable GTM alignment between product innovation and commercial results."
$new2 = "Transformed CRM strategy across financial services, hospitality, and retail portfolios by replacing legacy journey-based models with a dynamic, moments-based engagement framework—translating AI and predictive personalization capabilities into business value narratives that resonated with enterprise buyers. Architected automated content systems and GenAI-enabled tools that accelerated execution, improved message relevance, and delivered personalized brand experiences at scale. Drove measurable outcomes: 5% lift in customer retention, 12% increase in cross-sell performance, and 32% reduction in production costs."
Replace-ParagraphText $d $old2 $new2 | Out-Null

# Change 4
$old3 = "Directed the first global brand strategy platform and corporate campaign for a leading pharmaceutical company to raise awareness of an emerging Cell & Gene Therapy treatment class. Translated complex science into compelling executive narratives by partnering with R&D, corporate affairs, and biopharma leadership to develop a clear, human-centered strategic storytelling framework aligned with enterprise brand positioning. Launched an always-on integrated campaign supported by strategic benefit partnerships, driving GTM alignment and stakeholder engagement across policy, provider, and patient communities."
$new3 = "Created the first global brand strategy platform and corporate campaign for a leading pharmaceutical company, translating complex Cell & Gene Therapy science into a clear, human-centered narrative through partnership with R&D, corporate affairs, and biopharma leadership. Delivered an always-on integrated campaign with strategic benefit partnerships, driving GTM alignment and stakeholder engagement across policy, provider, and patient communities."
Replace-ParagraphText $d $old3 $new3 | Out-Null

# Change 5
$old4 = "Led the development of Deloitte's 'Only See Possible' brand platform, architecting a modular campaign ecosystem designed to unify executive narratives across industries and offerings. Directed strategic storytelling, narrative development, and executional rollout in partnership with internal creative, media, and performance teams, driving GTM alignment across functions. The platform delivered a 43% increase in lead generation and a 33% lift in brand relevancy across five priority sectors."
$new4 = "Led the development of Deloitte's 'Only See Possible' brand platform, architecting a modular campaign ecosystem designed to unify GTM alignment across industries and offerings. Directed strategic storytelling, narrative architecture, and executional rollout in partnership with internal creative, media, and performance teams. The platform drove a 43% increase in lead generation and a 33% lift in brand relevancy across five priority sectors."
Replace-ParagraphText $d $old4 $new4 | Out-Null

# Change 6
$old5 = "Created unified messaging architecture for B2B technology platform brand through strategic storytelling and creative strategy, delivering 15% awareness increase and 12% YoY acquisition growth"
$new5 = "Architected strategic storytelling framework and messaging architecture to unify B2B technology platform brand, driving GTM alignment across creative strategy that delivered 15% awareness increase and 12% YoY acquisition growth"
Replace-ParagraphText $d $old5 $new5 | Out-Null

# Change 7
$old6 = "Created a data-driven brand architecture system to unify Deloitte Consulting's executive narratives and positioning across its fastest-growing industry sectors. Designed the framework to ensure consistency of strategic storytelling, value proposition, and competitive positioning at both the enterprise and solution levels. The system enabled GTM alignment across verticals and contributed to a 3x increase in brand-led pipeline opportunities."
$new6 = "Architected a data-driven brand architecture system to unify Deloitte Consulting's competitive positioning across its fastest-growing industry sectors, designing the framework to ensure narrative consistency, value proposition clarity, and competitive differentiation at both enterprise and solution levels—enabling GTM alignment across verticals and contributing to a 3x increase in brand-led pipeline opportunities."
Replace-ParagraphText $d $old6 $new6 | Out-Null

# Change 8
$old7 = "Led strategic storytelling initiative for national dental services provider navigating category modernization, building segmentation model to identify high-potential patient segments and developing the 'Yes' platform—executive narratives brought to life through real patient stories that translated organizational values into compelling proof points around care, access, and experience. Activated GTM campaign across national and local TV, radio, paid search, and social, driving a 16% increase in demand and a 14% lift in revenue within the first year."
$new7 = "Led strategic storytelling initiative and creative platform development for a national dental services provider seeking to modernize its image and reestablish relevance in a competitive category. Created a segmentation model to identify high-potential patient segments and directed the `"Yes`" platform—translating brand purpose into executive narratives through real patient stories that highlighted the organization's culture of saying `"yes`" to care, access, and experience. Activated the campaign across national and local TV, radio, paid search, and social, driving a 16% increase in demand and a 14% lift in revenue within the first year."
Replace-ParagraphText $d $old7 $new7 | Out-Null
